# fix serious bug parsing date and money
# Column D ("時間") was incorrectly parsed as a numeric 0, while the
# real date-range text ended up being stuffed (scaled) into column E
# ("金額"). Restore the correct date-range text in column D and reset
# column E back to 0 for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> correct date-range text for column D.
# Column E is reset to 0 for every one of these rows.
$fixes = @{
    76  = "8-10月"
    77  = "8-10月"
    79  = "10-12月"
    81  = "6-8月；10-12月"
    83  = "9-11月"
    84  = "9-11月"
    86  = "9月；11-12月"
    89  = "10~12月"
    91  = "11月-103年1月"
    92  = "11月-103年1月"
    93  = "11月-103年1月"
    95  = "11月"
    96  = "11月"
    97  = "11月"
    99  = "12月"
    100 = "12月"
    102 = "12月"
    108 = "11月-103年1月"
    109 = "11月-103年1月"
    110 = "11月-103年1月"
    116 = "10-12月"
    117 = "10-12月"
}

foreach ($row in $fixes.Keys) {
    $ws.Range("D$row").Value = $fixes[$row]
    $ws.Range("E$row").Value = 0
}
